$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for rows 2-7, columns E (5) through T (20)
# Row layout: E,F,G,H,I,J,K,L,M,N,O,P,Q,R,S,T

$data = @{
    2 = @(3,1,4.278141666666667,12.834425,0.9663225094340192,0.9663225094340191,3,1,7.940402333333334,23.821207,0.1931648990487216,0.1931648990487216,33.97016607233056,305.731494650975,0.1866595899833297,0.1866595899833297)
    3 = @(3,1,4.278141666666667,12.834425,0.9663225094340192,0.9663225094340191,3,1,26.95384733333333,80.861542,0.6557019380820612,0.6557019380820612,115.3123773537055,1037.81139618335,0.6336195422482073,0.6336195422482073)
    4 = @(3,1,4.278141666666667,12.834425,0.9663225094340192,0.9663225094340191,3,1,6.212609666666666,18.637829,0.1511331628692172,0.1511331628692172,26.57842427370278,239.205818463325,0.1460433772024823,0.1460433772024823)
    5 = @(3,1,0.1490983333333333,0.447295,0.03367749056598092,0.03367749056598091,3,1,7.940402333333334,23.821207,0.1931648990487216,0.1931648990487216,1.183900753896111,10.655106785065,0.00650530906539198,0.006505309065391978)
    6 = @(3,1,0.1490983333333333,0.447295,0.03367749056598092,0.03367749056598091,3,1,26.95384733333333,80.861542,0.6557019380820612,0.6557019380820612,4.018773714321111,36.16896342889,0.02208239583385402,0.02208239583385402)
    7 = @(3,1,0.1490983333333333,0.447295,0.03367749056598092,0.03367749056598091,3,1,6.212609666666666,18.637829,0.1511331628692172,0.1511331628692172,0.9262897469505555,8.336607722555,0.005089785666734918,0.005089785666734918)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    $col = 5
    foreach ($v in $vals) {
        $ws.Cells.Item($r, $col).Value = $v
        $col = $col + 1
    }
}
